$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2525773195876289
$ws.Range("C2").Value = 0.4381443298969072
$ws.Range("J2").Value = 0.03092783505154639
$ws.Range("P2").Value = 0.1675257731958763
$ws.Range("S2").Value = 0.1108247422680412
$ws.Range("B3").Value = 0.005494505494505495
$ws.Range("C3").Value = 0.02747252747252747
$ws.Range("J3").Value = 0.02197802197802198
$ws.Range("P3").Value = 0.7747252747252747
$ws.Range("S3").Value = 0.1703296703296703
$ws.Range("J4").Value = 0.1296296296296296
$ws.Range("P4").Value = 0.5925925925925926
$ws.Range("S4").Value = 0.2777777777777778
$ws.Range("B6").Value = 0.05905511811023622
$ws.Range("D6").Value = 0.01181102362204724
$ws.Range("F6").Value = 0.06692913385826772
$ws.Range("J6").Value = 0.2637795275590551
$ws.Range("O6").Value = 0.02362204724409449
$ws.Range("Q6").Value = 0.1141732283464567
$ws.Range("R6").Value = 0.07874015748031496
$ws.Range("S6").Value = 0.3818897637795275
$ws.Range("B7").Value = 0.1179039301310044
$ws.Range("D7").Value = 0.02620087336244541
$ws.Range("F7").Value = 0.07860262008733625
$ws.Range("J7").Value = 0.1266375545851528
$ws.Range("O7").Value = 0.03056768558951965
$ws.Range("Q7").Value = 0.1222707423580786
$ws.Range("R7").Value = 0.09170305676855896
$ws.Range("S7").Value = 0.4061135371179039
$ws.Range("B8").Value = 0.1037735849056604
$ws.Range("D8").Value = 0.01886792452830189
$ws.Range("E8").Value = 0.001886792452830189
$ws.Range("F8").Value = 0.0660377358490566
$ws.Range("J8").Value = 0.09056603773584905
$ws.Range("O8").Value = 0.01320754716981132
$ws.Range("Q8").Value = 0.1716981132075472
$ws.Range("R8").Value = 0.1226415094339623
$ws.Range("S8").Value = 0.4113207547169812
$ws.Range("B9").Value = 0.09821428571428571
$ws.Range("D9").Value = 0.03125
$ws.Range("F9").Value = 0.0625
$ws.Range("J9").Value = 0.1383928571428572
$ws.Range("O9").Value = 0.04017857142857143
$ws.Range("Q9").Value = 0.1160714285714286
$ws.Range("R9").Value = 0.08482142857142858
$ws.Range("S9").Value = 0.4285714285714285
$ws.Range("B10").Value = 0.1233480176211454
$ws.Range("D10").Value = 0.02202643171806168
$ws.Range("E10").Value = 0.0007342143906020558
$ws.Range("F10").Value = 0.07195301027900147
$ws.Range("J10").Value = 0.1071953010279001
$ws.Range("O10").Value = 0.02055800293685756
$ws.Range("Q10").Value = 0.2320117474302496
$ws.Range("R10").Value = 0.06681350954478708
$ws.Range("S10").Value = 0.355359765051395
$ws.Range("G11").Value = 0.1534391534391534
$ws.Range("J11").Value = 0.07407407407407407
$ws.Range("K11").Value = 0.2195767195767196
$ws.Range("L11").Value = 0.5264550264550265
$ws.Range("S11").Value = 0.02645502645502645
$ws.Range("G12").Value = 0.7241379310344828
$ws.Range("J12").Value = 0.2413793103448276
$ws.Range("K12").Value = 0.004926108374384237
$ws.Range("L12").Value = 0.009852216748768473
$ws.Range("S12").Value = 0.01970443349753695
$ws.Range("F13").Value = 0.02040816326530612
$ws.Range("G13").Value = 0.5714285714285714
$ws.Range("J13").Value = 0.3265306122448979
$ws.Range("S13").Value = 0.08163265306122448
$ws.Range("F15").Value = 0.02262443438914027
$ws.Range("H15").Value = 0.1809954751131222
$ws.Range("I15").Value = 0.06334841628959276
$ws.Range("J15").Value = 0.3076923076923077
$ws.Range("K15").Value = 0.05429864253393665
$ws.Range("M15").Value = 0.009049773755656109
$ws.Range("N15").Value = 0.004524886877828055
$ws.Range("O15").Value = 0.05429864253393665
$ws.Range("S15").Value = 0.3031674208144796
$ws.Range("F16").Value = 0.008771929824561403
$ws.Range("H16").Value = 0.206140350877193
$ws.Range("I16").Value = 0.09649122807017543
$ws.Range("J16").Value = 0.3903508771929824
$ws.Range("K16").Value = 0.1008771929824561
$ws.Range("M16").Value = 0.01754385964912281
$ws.Range("O16").Value = 0.03070175438596491
$ws.Range("S16").Value = 0.1491228070175439
$ws.Range("F17").Value = 0.01844262295081967
$ws.Range("H17").Value = 0.1516393442622951
$ws.Range("I17").Value = 0.1004098360655738
$ws.Range("J17").Value = 0.4180327868852459
$ws.Range("K17").Value = 0.1086065573770492
$ws.Range("M17").Value = 0.02459016393442623
$ws.Range("N17").Value = 0.004098360655737705
$ws.Range("O17").Value = 0.04508196721311476
$ws.Range("S17").Value = 0.1290983606557377
$ws.Range("F18").Value = 0.0186046511627907
$ws.Range("H18").Value = 0.2418604651162791
$ws.Range("I18").Value = 0.07906976744186046
$ws.Range("J18").Value = 0.3953488372093023
$ws.Range("K18").Value = 0.1255813953488372
$ws.Range("M18").Value = 0.01395348837209302
$ws.Range("O18").Value = 0.05116279069767442
$ws.Range("S18").Value = 0.07441860465116279
$ws.Range("F19").Value = 0.01169993117687543
$ws.Range("H19").Value = 0.2188575361321404
$ws.Range("I19").Value = 0.08671713695801789
$ws.Range("J19").Value = 0.3448038540949759
$ws.Range("K19").Value = 0.1197522367515485
$ws.Range("M19").Value = 0.02202339986235375
$ws.Range("N19").Value = 0.0006882312456985547
$ws.Range("O19").Value = 0.05987611837577426
$ws.Range("S19").Value = 0.1355815554026153
